$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1028.3334
$ws.Range("I41").Value = 1313.4445
$ws.Range("J41").Value = 743.2222
$ws.Range("K41").Value = 1313.4445
$ws.Range("L41").Value = 743.2222
$ws.Range("M41").Value = -873.4445000000001
$ws.Range("N41").Value = -1623.2222
$ws.Range("H94").Value = 1000000000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1000000000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1000000000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1000000902
$ws.Range("H113").Value = 144915
$ws.Range("I113").Value = 335301.66
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 335301.66
$ws.Range("L113").Value = 2125
$ws.Range("M113").Value = -332047.66
$ws.Range("N113").Value = -8633
$ws.Range("H115").Value = 794.8333
$ws.Range("I115").Value = 742.25
$ws.Range("K115").Value = 2226.75
$ws.Range("M115").Value = -659.75
$ws.Range("H137").Value = 3031.2307
$ws.Range("I137").Value = 1875.375
$ws.Range("J137").Value = 4880.6
$ws.Range("K137").Value = 5626.125
$ws.Range("L137").Value = 14641.8
$ws.Range("M137").Value = -3076.125
$ws.Range("N137").Value = -19741.8
$ws.Range("H138").Value = 4396.271
$ws.Range("I138").Value = 1445.2963
$ws.Range("J138").Value = 8190.381
$ws.Range("K138").Value = 4335.8889
$ws.Range("L138").Value = 24571.143
$ws.Range("M138").Value = 804.1111000000001
$ws.Range("N138").Value = -34851.143
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3907.5
$ws.Range("I45").Value = 3297.818
$ws.Range("J45").Value = 4865.5713
$ws.Range("K45").Value = 3297.818
$ws.Range("L45").Value = 4865.5713
$ws.Range("M45").Value = -2920.818
$ws.Range("N45").Value = -5619.5713
$ws.Range("H74").Value = 3941.25
$ws.Range("I74").Value = 1950
$ws.Range("J74").Value = 5932.5
$ws.Range("K74").Value = 1950
$ws.Range("L74").Value = 5932.5
$ws.Range("M74").Value = -1076
$ws.Range("N74").Value = -7680.5
$ws.Range("H77").Value = 3941.25
$ws.Range("I77").Value = 1950
$ws.Range("J77").Value = 5932.5
$ws.Range("K77").Value = 9750
$ws.Range("L77").Value = 29662.5
$ws.Range("M77").Value = -5382
$ws.Range("N77").Value = -38398.5
$ws.Range("H122").Value = 3354.8235
$ws.Range("I122").Value = 3815.6667
$ws.Range("J122").Value = 3103.4546
$ws.Range("K122").Value = 11447.0001
$ws.Range("L122").Value = 9310.363799999999
$ws.Range("M122").Value = -8997.000100000001
$ws.Range("N122").Value = -14210.3638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 34716.613
$ws.Range("I20").Value = 48825.047
$ws.Range("K20").Value = 48825.047
$ws.Range("M20").Value = -48578.047
$ws.Range("H64").Value = 366.1111
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 366.1111
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 366.1111
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -816.1111000000001
$ws.Range("H67").Value = 366.1111
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 366.1111
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 366.1111
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -1926.1111
$ws.Range("H107").Value = 22227320
$ws.Range("I107").Value = 47623020
$ws.Range("J107").Value = 6083.375
$ws.Range("K107").Value = 47623020
$ws.Range("L107").Value = 6083.375
$ws.Range("M107").Value = -47621100
$ws.Range("N107").Value = -9923.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4300
$ws.Range("I86").Value = 4200
$ws.Range("J86").Value = 4350
$ws.Range("K86").Value = 4200
$ws.Range("L86").Value = 4350
$ws.Range("M86").Value = -3077
$ws.Range("N86").Value = -6596
$ws.Range("H89").Value = 4300
$ws.Range("I89").Value = 4200
$ws.Range("J89").Value = 4350
$ws.Range("K89").Value = 21000
$ws.Range("L89").Value = 21750
$ws.Range("M89").Value = -15384
$ws.Range("N89").Value = -32982
$ws.Range("H99").Value = 20125.5
$ws.Range("I99").Value = 2926.6667
$ws.Range("J99").Value = 37324.332
$ws.Range("K99").Value = 2926.6667
$ws.Range("L99").Value = 37324.332
$ws.Range("M99").Value = -1428.6667
$ws.Range("N99").Value = -40320.332
$ws.Range("H122").Value = 3234.724
$ws.Range("I122").Value = 3099.7
$ws.Range("K122").Value = 9299.099999999999
$ws.Range("M122").Value = -6849.099999999999
$ws.Range("H126").Value = 20125.5
$ws.Range("I126").Value = 2926.6667
$ws.Range("J126").Value = 37324.332
$ws.Range("K126").Value = 8780.000100000001
$ws.Range("L126").Value = 111972.996
$ws.Range("M126").Value = -6310.000100000001
$ws.Range("N126").Value = -116912.996
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 828.23
$ws.Range("I131").Value = 539.75
$ws.Range("J131").Value = 840.25
$ws.Range("K131").Value = 1619.25
$ws.Range("L131").Value = 2520.75
$ws.Range("M131").Value = 3420.75
$ws.Range("N131").Value = -12600.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2512.9565
$ws.Range("I126").Value = 2717
$ws.Range("J126").Value = 2290.3635
$ws.Range("K126").Value = 8151
$ws.Range("L126").Value = 6871.0905
$ws.Range("M126").Value = -5681
$ws.Range("N126").Value = -11811.0905
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2087.0476
$ws.Range("I7").Value = 1385
$ws.Range("J7").Value = 2725.2727
$ws.Range("K7").Value = 1385
$ws.Range("L7").Value = 2725.2727
$ws.Range("M7").Value = -1273
$ws.Range("N7").Value = -2949.2727
$ws.Range("H40").Value = 73720.57000000001
$ws.Range("I40").Value = 251772.5
$ws.Range("J40").Value = 2499.8
$ws.Range("K40").Value = 251772.5
$ws.Range("L40").Value = 2499.8
$ws.Range("M40").Value = -251636.5
$ws.Range("N40").Value = -2771.8
$ws.Range("H46").Value = 1125734.5
$ws.Range("I46").Value = 920
$ws.Range("J46").Value = 2531752.5
$ws.Range("K46").Value = 920
$ws.Range("L46").Value = 2531752.5
$ws.Range("M46").Value = -732
$ws.Range("N46").Value = -2532128.5
$ws.Range("H126").Value = 2087.0476
$ws.Range("I126").Value = 1385
$ws.Range("J126").Value = 2725.2727
$ws.Range("K126").Value = 4155
$ws.Range("L126").Value = 8175.8181
$ws.Range("M126").Value = -1685
$ws.Range("N126").Value = -13115.8181
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1981.0769
$ws.Range("I126").Value = 2596.6667
$ws.Range("J126").Value = 1453.4286
$ws.Range("K126").Value = 7790.000100000001
$ws.Range("L126").Value = 4360.2858
$ws.Range("M126").Value = -5320.000100000001
$ws.Range("N126").Value = -9300.2858
